$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("H2").Value = 3.25
$ws.Range("K2").Value = 2
$ws.Range("N2").Value = 7.5
$ws.Range("Q2").Value = 2.3
$ws.Range("R2").Value = 1.6
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 1.75
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 8.5
$ws.Range("AA2").Value = 19
$ws.Range("AB2").Value = 34
$ws.Range("AC2").Value = 7.5
$ws.Range("AD2").Value = 6
$ws.Range("AF2").Value = 67
$ws.Range("AG2").Value = 401
$ws.Range("AH2").Value = 9.5
$ws.Range("AJ2").Value = 15
$ws.Range("AT2").Value = 2.5
$ws.Range("AU2").Value = 9
$ws.Range("BA2").Value = 126
$ws.Range("BB2").Value = 301

# Row 3 updates
$ws.Range("L3").Value = 6
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("AC3").Value = 8
$ws.Range("AF3").Value = 67
$ws.Range("AV3").Value = 67
$ws.Range("AZ3").Value = 126
$ws.Range("BA3").Value = 151
